# Apply updated cryptocurrency data to the worksheet.
# Column D (Price) and column E (Volume/1h) values are textual (inline strings)
# in the source data, and some of the Price values look like plain numbers
# (e.g. "1.00", "6.40") or multi-dot grouped numbers (e.g. "63.349.58").
# To prevent Excel's automatic type coercion from turning these into numeric
# values (which would strip formatting like trailing zeros), the number
# format of columns D and E is forced to Text ("@") before the values are
# written, and then reset back to the default "Normal" style afterwards so
# the cell styling in the file is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "63.349.58"
    "E2" = "  +1.66%  "
    "D3" = "3.470.14"
    "E3" = "  +1.11%  "
    "D4" = "0.999"
    "E4" = "  +0.02%  "
    "D5" = "581.69"
    "E5" = "  +0.27%  "
    "D6" = "147.84"
    "E6" = "  +1.77%  "
    "D7" = "3.469.77"
    "E7" = "  +1.08%  "
    "E8" = "  -0.03%  "
    "E9" = "  +0.37%  "
    "D10" = "7.74"
    "E10" = "  +1.65%  "
    "E11" = "  +0.64%  "
    "E12" = "  +4.05%  "
    "D13" = "4.066.16"
    "E13" = "  +1.19%  "
    "D14" = "29.56"
    "E14" = "  +1.86%  "
    "E15" = "  +2.65%  "
    "D16" = "3.474.72"
    "E16" = "  +1.31%  "
    "E17" = "  +0.52%  "
    "D18" = "63.307.58"
    "E18" = "  +1.61%  "
    "D19" = "6.40"
    "E19" = "  +2.91%  "
    "E20" = "  +3.14%  "
    "D21" = "9.34"
    "E21" = "  +1.18%  "
    "D22" = "389.27"
    "E22" = "  -1.38%  "
    "E23" = "  +1.84%  "
    "D24" = "74.50"
    "E24" = "  -0.64%  "
    "E25" = "  -0.14%  "
    "D26" = "3.614.68"
    "E26" = "  +1.30%  "
    "E27" = "  -0.18%  "
    "E28" = "  -2.80%  "
    "D29" = "7.66"
    "E29" = "  +1.34%  "
    "D30" = "1.00"
    "E30" = "  +0.24%  "
    "E31" = "  +2.01%  "
    "E32" = "  -0.44%  "
    "E33" = "  +0.01%  "
    "E34" = "  -4.24%  "
    "D35" = "23.49"
    "E35" = "  -0.75%  "
    "B36" = "ImmutableX"
    "C36" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D36" = "1.62"
    "E36" = "  +6.26%  "
    "B37" = "NEARProtocol"
    "C37" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D37" = "5.33"
    "E37" = "  -0.17%  "
    "B38" = "Aptos"
    "C38" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D38" = "7.16"
    "E38" = "  +1.97%  "
    "D39" = "32.02"
    "E39" = "  +10.63%  "
    "D40" = "168.62"
    "E40" = "  +0.54%  "
    "D41" = "3.509.53"
    "E41" = "  +1.30%  "
    "D42" = "0.0766"
    "E42" = "  +1.30%  "
    "D43" = "0.794"
    "E43" = "  +0.57%  "
    "E44" = "  +3.55%  "
    "D45" = "42.42"
    "E46" = "  +2.81%  "
    "E47" = "  -1.49%  "
    "D48" = "2.595.62"
    "E48" = "  +3.16%  "
    "E49" = "  +8.14%  "
    "E50" = "  +1.95%  "
    "D51" = "23.03"
    "E51" = "  -0.57%  "
}

# Force text format on the affected columns so numeric-looking strings are
# preserved exactly as given (no silent conversion to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# Restore the default cell style so no stray number-format/style is left
# behind on cells that didn't have one originally.
$ws.Range("D2:E51").Style = "Normal"
